# Update weekly consolidated fruit/vegetable price data
# (Hortaliza, Agrícola del Norte S.A. de Arica - Ramas de apio)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => Fecha(D), Variedad(H), Calidad(I), Volumen(J), Precio min(K), Precio max(L), Precio prom(M), Precio $/Kg(P)
$rows = @(
    @{ Row = 2;  D = 44575; H = "Sin especificar"; I = "Primera"; J = 160; K = 6500; L = 7000; M = 6750; P = 6750 }
    @{ Row = 3;  D = 44497; H = "Sin especificar"; I = "Primera"; J = 160; K = 5000; L = 6000; M = 5500; P = 5500 }
    @{ Row = 4;  D = 44263; H = "Sin especificar"; I = "Primera"; J = 100; K = 7000; L = 8000; M = 7500; P = 7500 }
    @{ Row = 6;  D = 44371; H = "Sin especificar"; I = "Primera"; J = 80;  K = 7000; L = 8000; M = 7375; P = 7375 }
    @{ Row = 7;  D = 44636; H = "Americana (o)";   I = "Primera"; J = 60;  K = 8000; L = 9000; M = 8500; P = 8500 }
    @{ Row = 8;  D = 44559; H = "Americana (o)";   I = "Primera"; J = 100; K = 5000; L = 6000; M = 5500; P = 5500 }
    @{ Row = 9;  D = 44414; H = "Sin especificar"; I = "Primera"; J = 100; K = 6000; L = 7000; M = 6500; P = 6500 }
    @{ Row = 10; D = 44309; H = "Sin especificar"; I = "Primera"; J = 50;  K = 8000; L = 9000; M = 8500; P = 8500 }
    @{ Row = 11; D = 44410; H = "Sin especificar"; I = "Primera"; J = 100; K = 5500; L = 6000; M = 5750; P = 5750 }
    @{ Row = 12; D = 44253; H = "Americana (o)";   I = "Segunda"; J = 100; K = 4000; L = 4500; M = 4250; P = 4250 }
    @{ Row = 13; D = 44699; H = "Sin especificar"; I = "Primera"; J = 50;  K = 9000; L = 9500; M = 9250; P = 9250 }
    @{ Row = 14; D = 44281; H = "Sin especificar"; I = "Primera"; J = 100; K = 5000; L = 6000; M = 5500; P = 5500 }
    @{ Row = 15; D = 44259; H = "Sin especificar"; I = "Primera"; J = 80;  K = 4000; L = 4500; M = 4250; P = 4250 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value  = $r.D   # D - Fecha
    $ws.Cells.Item($row, 8).Value  = $r.H   # H - Variedad
    $ws.Cells.Item($row, 9).Value  = $r.I   # I - Calidad
    $ws.Cells.Item($row, 10).Value = $r.J   # J - Volumen
    $ws.Cells.Item($row, 11).Value = $r.K   # K - Precio mínimo
    $ws.Cells.Item($row, 12).Value = $r.L   # L - Precio máximo
    $ws.Cells.Item($row, 13).Value = $r.M   # M - Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $r.P   # P - Precio $/Kg
}
